# Merge the separately-formatted "<id>", "<pageID>", "</id>" runs into a
# single run per occurrence, e.g. "<id>p049v_1</id>", keeping the
# Courier-New / color-7f6000 formatting of the opening "<id>" run.
#
# wdReplaceAll = 2 (word:Replace parameter)
# wdFindContinue = 1 (word:Wrap parameter)

$d = $word.ActiveDocument

$ids = @("p049v_1", "p049v_2", "p049v_3")

foreach ($pageId in $ids) {
    $find = "<id>" + $pageId + "</id>"
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $find, 2) | Out-Null
}
